# Update odds values on "Sheet1" to reflect the latest FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 (Aguilas - Dep. Cali)
$ws.Range("Q5").Value = 2.3
$ws.Range("R5").Value = 1.6

# Row 6 (Dorados de Sinaloa - Tlaxcala)
$ws.Range("G6").Value  = 2.15
$ws.Range("I6").Value  = 3.25
$ws.Range("J6").Value  = 2.7
$ws.Range("L6").Value  = 3.65
$ws.Range("Q6").Value  = 1.98
$ws.Range("W6").Value  = 7.3
$ws.Range("X6").Value  = 10.25
$ws.Range("Z6").Value  = 21
$ws.Range("AD6").Value = 6.1
$ws.Range("AG6").Value = 9.75
$ws.Range("AH6").Value = 17.5
$ws.Range("AI6").Value = 11.25
$ws.Range("AJ6").Value = 45
$ws.Range("AK6").Value = 29
$ws.Range("AQ6").Value = 40
$ws.Range("AT6").Value = 2.6
$ws.Range("AU6").Value = 6.6
$ws.Range("AW6").Value = 5.2
$ws.Range("AX6").Value = 17.5
$ws.Range("AY6").Value = 22
$ws.Range("BB6").Value = 250

# Row 12 (Academico Viseu - Maritimo)
$ws.Range("N12").Value = 10
